# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 - the "Status"/locale columns).
# 2) Narrow the Status column(s) to match the new (shorter) text:
#    Overview columns E & F, and column C on the zh-cn / de-de sheets.
#    Target stored width ~= 13.4101848602295 chars; the COM ColumnWidth
#    property is quantized to 1/6-character increments by this host, so
#    12.5 is the input that lands on the closest reachable stored width
#    (13.333333333333334).

$wb = $excel.ActiveWorkbook

$targetColumnWidth = 12.5

# --- Sheet "Overview": Status values live in columns E (zh-cn) & F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# --- Sheet "zh-cn": Status values live in column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth

# --- Sheet "de-de": Status values live in column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
